# Build site update for LOM3256.xlsx
# The syllabus content was trimmed down; this deletes the now-empty "Docentes
# responsaveis" value row (old row 13, which only carried the professor's name
# in columns B/C with no label in column A) and shifts everything below it up
# by one row. The remaining long-form paragraphs (objectives, short/long
# syllabus, bibliography, etc.) are also replaced with much shorter content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row that used to hold the professor's name underneath
# "Docentes responsaveis:" (it had no A-column label) -- this shifts rows
# 14-24 up to become rows 13-23.
$ws.Rows(13).Delete() | Out-Null

# --- Overwrite the cells whose text content changed ---

# "Objetivos:" (row 10) now just shows the responsible professor's id/name
# instead of the long paragraph describing the course objectives.
$ws.Range("B10").Value2 = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C10").Value2 = "1176388 - Luiz Tadeu Fernandes Eleno"

# "Programa resumido:" (now row 13) becomes just "Semestral".
$ws.Range("B13").Value2 = "Semestral"
$ws.Range("C13").Value2 = "Semestral"

# "Programa:" (now row 15) becomes the activation date.
$ws.Range("B15").Value2 = "15/07/2015"
$ws.Range("C15").Value2 = "15/07/2015"

# "Metodo:" (now row 18) becomes the responsible professor's id/name again.
$ws.Range("B18").Value2 = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value2 = "1176388 - Luiz Tadeu Fernandes Eleno"

# "Criterio:" (now row 19) takes what used to be under "Metodo:".
$ws.Range("B19").Value2 = "Aulas expositivas, trabalhos e exercícios comentados."
$ws.Range("C19").Value2 = "Aulas expositivas, trabalhos e exercícios comentados."

# "Norma de recuperacao:" (now row 20) takes what used to be under "Criterio:".
$ws.Range("B20").Value2 = "Média aritmética de trabalhos propostos ao longo do curso."
$ws.Range("C20").Value2 = "Média aritmética de trabalhos propostos ao longo do curso."

# "Bibliografia:" (now row 21) takes what used to be under "Norma de recuperacao:".
$ws.Range("B21").Value2 = "Não haverá exame de recuperação"
$ws.Range("C21").Value2 = "Não haverá exame de recuperação"
